$d = $word.ActiveDocument

# --- Edit 1: "Can cu Luat Quan ly thue ngay 13 thang 6 nam 2019;" ---
# Split into: "Can cu Luat Quan ly " + "<luat_qlt_ngay>" + ";"
# with a _GoBack bookmark wrapping the whole paragraph content.

$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Căn cứ Luật Quản lý thuế ngày 13 tháng 6 năm 2019;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # rng now covers the whole matched text (the single run).
    $paraStart = $rng.Start

    # Collapse to start, insert bookmark start marker there (we'll add bookmark after constructing runs)
    $rng.Text = "Căn cứ Luật Quản lý ;"

    # Now rng covers the new text; we need to insert the placeholder between "Căn cứ Luật Quản lý " and ";"
    $insPoint = $paraStart + [int]"Căn cứ Luật Quản lý ".Length
    $insRange = $d.Range($insPoint, $insPoint)
    $insRange.InsertAfter("<luat_qlt_ngay>")

    # Apply lang nl-NL formatting to the placeholder run
    $placeholderRange = $d.Range($insPoint, $insPoint + [int]"<luat_qlt_ngay>".Length)
    $placeholderRange.LanguageID = 19

    # Add bookmark spanning entire paragraph text (from paraStart to end of ";")
    $bmEnd = $insPoint + [int]"<luat_qlt_ngay>".Length + 1
    $bmRange = $d.Range($paraStart, $bmEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- Edit 2: remove bookmark split in "Thuc hien giam s" + "at doi voi..." ---
$d.Bookmarks("_GoBack").Delete()
